$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column P. This shifts the existing P..AA
#    block one column to the right (P->Q, Q->R, ... AA->AB) and extends the
#    sheet's used range out to column AB.
$ws.Range("P1").EntireColumn.Insert()

# 2. Header for the freshly inserted column.
$ws.Range("P1").Value = "Canola.FrostHeatGrainWt"

# New column takes on the same width as its neighbour (Canola.Grain.Wt).
$ws.Columns("P").ColumnWidth = $ws.Columns("O").ColumnWidth

# 3. The harvest rows (where Canola.Grain.Wt / column O is populated) get the
#    new Canola.FrostHeatGrainWt value seeded from the Grain.Wt value.
$ws.Range("P8").Value = $ws.Range("O8").Value2
$ws.Range("P14").Value = $ws.Range("O14").Value2
$ws.Range("P27").Value = $ws.Range("O27").Value2
$ws.Range("P28").Value = $ws.Range("O28").Value2
$ws.Range("P29").Value = $ws.Range("O29").Value2

# 4. Column H (Canola.DaysAfterSowing) no longer carries an explicit
#    (redundant General) number-format style.
$ws.Range("H2:H29").ClearFormats()
$ws.Range("H41").ClearFormats()

# 5. Re-establish the AutoFilter over the new, wider range.
$ws.AutoFilterMode = $false
$ws.Range("A1:AB29").AutoFilter() | Out-Null

# 6. The hidden _FilterDatabase defined name also needs to track the wider
#    range (Excel keeps this in lock-step with the AutoFilter, but it is not
#    automatic here).
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Observed!_FilterDatabase") {
    $n.RefersTo = "=Observed!`$A`$1:`$AB`$29"
  }
}

# 7. Page setup (printer/paper size) picked up from the default printer --
#    A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
